$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D3 (was a numeric value, becomes blank/empty)
$ws.Range("D3").Value = $null

# New row 8: carry row 7's formatting down, then fill in the
# corrected "Other" row (done before row 7 is renamed to "Biogas").
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 345.1642501212369

# Row 7: rename "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 35.44455893018135
